$d = $word.ActiveDocument

# The document currently ends with a ListParagraph bullet (numId=1):
#   "Frontend deployed successfully!"
# We append four new paragraphs after it:
#   1. "3/28/2024"                                              (plain paragraph)
#   2. "Working on deployment, using Render for backend."       (bullet)
#   3. "Trying to deploy rabbitmq -> Deployed successfully"     (bullet, w/ Wingdings arrow)
#   4. "Trying to connect rabbitmq with server"                 (bullet)
#
# New paragraphs inserted right after a ListParagraph/numId=1 paragraph
# automatically inherit that same list formatting, so we build all four
# while still chained off list paragraphs, and only afterwards (a) strip
# the list formatting back off paragraph #1 (the date line) and (b) mark
# the arrow character as Wingdings. Both of those are deferred to the end
# of the script because touching .Font / .Style on a range changes the
# "current typing formatting" and would otherwise bleed into whatever
# gets typed next.

# ---- 1. "3/28/2024" ----------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("3/28/2024")
$dateParaIndex = $d.Paragraphs.Count

# ---- 2. "Working on deployment, using Render for backend." -------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Working on deployment, using Render for backend.")

# ---- 3. "Trying to deploy rabbitmq" + Wingdings arrow + "Deployed successfully"
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Trying to deploy rabbitmq ")
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertAfter([char]0xF0E0)
$arrowParaIndex = $d.Paragraphs.Count
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertAfter(" Deployed successfully")

# ---- 4. "Trying to connect rabbitmq with server" ------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Trying to connect rabbitmq with server")

# ---- Cleanup pass (deferred so it can't bleed into later insertions) ---

# (a) Re-font the single arrow character as Wingdings.
$arrowPara = $d.Paragraphs.Item($arrowParaIndex)
$paraText = $arrowPara.Range.Text
$symOffsetInPara = $paraText.IndexOf([char]0xF0E0)
$paraStart = $arrowPara.Range.Start
$symStart = $paraStart + $symOffsetInPara
$symRange = $d.Range($symStart, $symStart + 1)
$symRange.Font.Name = "Wingdings"

# (b) Strip the date paragraph back to a plain (non-list) paragraph.
$datePara = $d.Paragraphs.Item($dateParaIndex)
$datePara.Range.ListFormat.RemoveNumbers()
$datePara.Style = $d.Styles("Normal")

Write-Output "done"
